$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.926.47"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "1.815.52"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'310.10"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4651"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").Value = "'0.3706"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").Value = "'0.07371"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'0.8737"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "1.846.48"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "'5.356"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "'6.521"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "'0.07061"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "'91.74"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'0.000008734"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "26.932.73"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "'5.316"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").Value = "2.044.33"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'1.908"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").Value = "'151.93"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'2.152"
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("D29").Value = "'5.330"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").Value = "'115.88"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").Value = "'0.08921"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'0.7585"
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").Value = "'1.156"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").Value = "'4.475"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'2.919"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "'0.05268"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'2.419"
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").Value = "'2.943"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5357"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'7.249"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'0.1666"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").Value = "'8.458"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "'10.37"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'103.33"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").Value = "'0.06293"
$ws.Range("E51").Value = "  -0.88%  "